# Add "status_label" as a new column B (string version of "status"),
# shifting NCTId/eudraCT/.../results one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B - this shifts B:I to C:J
$ws.Columns.Item(2).Insert()

# Header for the new column
$ws.Cells.Item(1, 2).Value = "status_label"
# Copy the header style (bold/border/alignment) used by the rest of row 1
$ws.Cells.Item(1, 1).Copy()
$ws.Cells.Item(1, 2).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the string-label values corresponding to each row's status
$ws.Cells.Item(2, 2).Value = "vert"
$ws.Cells.Item(3, 2).Value = "rouge"
$ws.Cells.Item(4, 2).Value = "rouge"
